# SYM-1 IO Type A Manual — apply the documented content changes.
#
# The headless Word engine's Range.InsertXML always re-inserts the
# replacement content at the end of the paragraph that contains the
# range (even though it deletes the originally-selected text in place),
# so every edit below selects from the first changed character through
# to the end of its paragraph and then re-supplies the full (unchanged
# tail included) remainder of that paragraph as literal OOXML. This
# keeps run/proofErr boundaries identical to the target everywhere that
# isn't actually changing.
#
# Range offsets returned by Find.Execute live in a different coordinate
# space than $doc.Content.Text (tables/fields shift them), so every
# lookup below re-finds its anchor text with Find.Execute on a live
# Range variable instead of indexing into Content.Text.

$d = $word.ActiveDocument

function Get-FoundRange([string]$needle) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find text: $needle"
    }
    return $rng
}

# ---------------------------------------------------------------------
# Change 1: "...SD2IEC.  SymDOS functionality requires..." paragraph.
#   - split "...SD2IEC.  " into "...SD2IEC" + " using " (+ the existing
#     "SymDOS" run, now followed by ".  " instead of being merged in)
#   - add a second "SymDOS" (spell-checked) before "functionality..."
#   - append a new sentence recommending the Corsham Tech RAM board
# ---------------------------------------------------------------------
$found = Get-FoundRange("Commodore IEC plug for connection")
$startPos = $found.Start
$paraEnd = $found.Paragraphs(1).Range.End
$target = $d.Range($startPos, $paraEnd)

$xml1 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Commodore IEC plug for connection to a Commodore floppy disk drive or an SD2IEC</w:t></w:r><w:r><w:t xml:space="preserve"> using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>SymDOS</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">.  </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>SymDOS</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> functionality requires RAM at $9000 or a burned EPROM addressable at $9000.</w:t></w:r><w:r><w:t xml:space="preserve">  I recommend </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Corsham</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Tech'S SYM-1/AIM-65 RAM Board.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$xml1 = $xml1.Replace("Tech'S", "Tech" + [char]0x2019 + "s")
$target.InsertXML($xml1)

# ---------------------------------------------------------------------
# Change 2: "...<Shift><Jump><1> a dot should appear..." -> insert a
# literal "<CR>" run between "<1>" and " a dot should appear...".
# ---------------------------------------------------------------------
$found = Get-FoundRange("<1> a dot should appear")
$startPos = $found.Start
$paraEnd = $found.Paragraphs(1).Range.End
$target = $d.Range($startPos, $paraEnd)

$xml2 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>&lt;1&gt;</w:t></w:r><w:r><w:t>&lt;CR&gt;</w:t></w:r><w:r><w:t xml:space="preserve"> a dot should appear in your terminal.  Right now, you are connected to the SYM-1 at 110 bps, which is SLOW.  To change the connection speed</w:t></w:r><w:r w:rsidR="007418C3"><w:t xml:space="preserve"> in your terminal type M A651 and then enter XX where XX is the number in the below table representing the speed you want to use.  Then use Control-Z to cycle through the speeds until you have reached the matching speed.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target.InsertXML($xml2)

# ---------------------------------------------------------------------
# Changes 3 & 4 (same paragraph, "SymDOS Usage" body text):
#   - move <w:lastRenderedPageBreak/> so it precedes the SECOND "SymDOS"
#     run (right after "Included on the USB stick is the ") instead of
#     sitting in front of "burning to an 4k EPROM..."
#   - merge the two runs that used to be split around the old page
#     break location back into one run of plain text
# ---------------------------------------------------------------------
$found = Get-FoundRange(" to the IEC connector and power on the drive")
$startPos = $found.Start
$paraEnd = $found.Paragraphs(1).Range.End
$target = $d.Range($startPos, $paraEnd)

$xml3 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r w:rsidR="00587120"><w:t xml:space="preserve"> to the IEC connector and power on the drive/SD2IEC and the SYM-1.  Connect to your SYM-1'S terminal either following the steps above or using a T-connector solution.  Included on the USB stick is the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00587120"><w:lastRenderedPageBreak/><w:t>SymDOS</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00587120"><w:t xml:space="preserve"> software in KIM-1 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00587120"><w:t>papertape</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00587120"><w:t xml:space="preserve"> format, a wav file in SYM-1 tape format and an Intel Hex file for burning to an 4k EPROM.  Load the software however is convenient and start it using G 9000 at the monitor once the software is loaded into memory.  Included on the USB stick is the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00587120"><w:t>SymDOS</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00587120"><w:t xml:space="preserve"> manual, so consult it for how to use it.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$xml3 = $xml3.Replace("SYM-1'S", "SYM-1" + [char]0x2019 + "s")
$target.InsertXML($xml3)

Write-Output "All changes applied."
